# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Ají" (Inferno, Primera,
# Región de Arica y Parinacota) above the current row 295, pushing the
# existing rows 295-308 down to 296-309.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(295).Insert()

$ws.Range("A295").Value = 9
$ws.Range("B295").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C295").Value = "Metropolitana"
$ws.Range("D295").Value = 44746
$ws.Range("E295").Value = 13
$ws.Range("F295").Value = 100112021
$ws.Range("G295").Value = "Ají"
$ws.Range("H295").Value = "Inferno"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 61
$ws.Range("K295").Value = 16000
$ws.Range("L295").Value = 18000
$ws.Range("M295").Value = 17016
$ws.Range("N295").Value = "`$/caja 12 kilos"
$ws.Range("O295").Value = "Región de Arica y Parinacota"
$ws.Range("P295").Value = 1418
$ws.Range("Q295").Value = 12
$ws.Range("R295").Value = "Hortaliza"
